$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typos in existing long-text cells ---
$c2text = @'
Die Technologie enttwickelt sichg stetig weiter. Ein wichtiger Aspekt der neuen Technologie ist die KI. Sieht wird in etlichen Bereichen eine Wende bewirken. Ein großer Aspekt sind Hier bei die Cybersecurity , Datenschutz, Big Data und die Infustrie 5.0 
'@
$ws.Range("C2").Value = $c2text

$c3text = @'
Die Gesellschaft muss inn Sachen Energie eine neue Richtung einschlagen. Es müssen neue Ideen geschaffen werden, um sauberer Energie  erzeugen.
'@
$ws.Range("C3").Value = $c3text

$c4text = @'
Arbeiter wollen eine gute Work Life Balance. Ihnen ist die zeitliche Vereinbarung zwischen privat und Arbeit wichtiger geworden. Dazu entwickeln Unternehmen flexible Arbeitsmodelle. Die weiter enteckelte Technologie und KI helgen bei dem Vorhaben. Jedoch gibt es einige negative Dinge wie die Privatsspähre, Datenschutz und Cybe Security.
'@
$ws.Range("C4").Value = $c4text

# --- Fill Maxi01/02/03 rows (6,7,8) ---
$b6text = @'
Umweltresilientere Stadtplanung
'@
$ws.Range("B6").Value = $b6text
$c6text = @'
Der Klimawandel schreitet stetig voran. Überflutungen, extreme Hitze und schlechter werdende Luftqualität durch Waldbrände und Smog sind eine Große Herausforderung für die Städte der Zukunft. Dabei sind nachhaltigere Technologien und Bauweisen immer mehr gefragt. Der Bau von vertical gardens an Außen und Innenwänden sorgt für ein angenehmeres Klima inmitten der Städte. Zugleich werden durch fortschrittlicherer Gentechnik, Pflanzen gezüchtet und in den vertical gardens eingesetzt, die mehr CO2 binden können und somit dem Klimawandel zusätzlich entgegenwirken soll.
'@
$ws.Range("C6").Value = $c6text

$b7text = @'
Alte Leute - es werden immer mehr
'@
$ws.Range("B7").Value = $b7text
$c7text = @'
Durch immer bahnbrechndere Erkenntnisse und Erfindungen im Bereich der Biotechnologie und Gentechnik kann die Menschheit ihren Todestag immer weiter nach hinten verschieben. Altersbedingte Krankheiten können leichter geheilt/ behandelt werden. Diese Entwicklung hat allerdings den negativen Effekt der Überalterung der Gesellschaft. Der Demografische Wandel wird verstärkt und etablierte Systeme, wie die Rente für alle werden zunehmend auf die Probe gestellt. Demzufolge steigt auch die Altersarmut weiter an und sorgt somit für eine weitere Belastung der Sozial und Gesundheitssysteme.
'@
$ws.Range("C7").Value = $c7text

$b8text = @'
Pflanzen - Grundpfeiler der Ernährung
'@
$ws.Range("B8").Value = $b8text
$c8text = @'
Weniger Biodiversität und mehr Monokulturen, das ist bislang das Gebot der Stunde. Allerdings führt das gerade in der Landwirtschaft zu einem Artensterben von Insekten, die wie die Bienen Pflanzen untereinander bestäuben und somit für eine Reproduktion dieser führt. Gentechnik wird dabei die Aufgabe der Bienen übernehmen müssen um weiterhin reichhaltige Ernten zu erhalten. Dabei rückt auch die Produktion von Pflanzen in den Mittelpunkt, die für den Verzehr durch den Menschen genutzt werden, statt für Nutztiere. Um den CO2 Footprint in der Landwirtschaft nachhaltig zu senken werden Nutztiere immer weiter durch Pflanziche oder im Labor gezüchtete Präperate ersetzt. Das führt zu einer Klima Entlastung im Bereich Wasserwirtschaft und Landwirtschaft. Positiver Nebeneffekt ist eine gesündere Bevölkerung, da Herzkreislauferkrankungen und Krebsentwicklungen seltener auftreten durch eine fleischlosere Ernährung.
'@
$ws.Range("C8").Value = $c8text

# --- Fill Jan R01/02/03 rows (18,19,20) ---
$b18text = @'
Digitale Transformation
'@
$ws.Range("B18").Value = $b18text
$c18text = @'
In dieser Zukunft verschmelzen KI und Web3.0 nahtlos, während die Menschheit der Zukunft auf KI-gestützte Bildung und "wearable technologies" setzt. Smart Homes, die über 5G und IoT vernetzt sind, optimieren Ressourcennutzung und tragen zur Energieeffizienz bei. Blockchain und eine ausgereifte Cybersecurity sichern persönliche Daten und fördern digitales Vertrauen. Klimabewusstsein treibt die F&E voran, insbesondere in Smart Cities, wo AR einen bedeutenden Beitrag zum Fortschritt leistet. Noch nie war die Menschheit so eng miteinander vernetzt. Doch im Schatten des Fortschritts leiden Menschen unter digitaler Erschöpfung. Autonome Fahrzeuge und smarte Technologien sind Teil des Alltags, stellen aber auch ethische Fragen in einer Welt, die sich rasend schnell digitalisiert.
'@
$ws.Range("C18").Value = $c18text

$b19text = @'
Wandel der Ökonomie und Lebensräume
'@
$ws.Range("B19").Value = $b19text
$c19text = @'
In der Ära der Industriellen Transformation formt eine durchdachte Infrastruktur und nachhaltige Technologien das allgemeine Leben. Stadtwerke führen die Energiewende an, gestützt durch flexible Energiespeicher und fortschrittliche Wasserstofftechnologien. Elektromobilität wird zum nahtlosen Service, angetrieben von klimaneutraler Solarenergie und Wasserstoff. Regulatorische Maßnahmen und innovativer Wissensaustausch sind  Schlüssel in diesem Paradigmenwechsel. Während durch fortgeschrittene Technologie CO2 zur Energiegewinnung genutzt wird und die globalen CO2-Emissionen sinken, werden Wasserressourcen neuartig aufbereitet um die Kreislaufiwrtschaft zu stärken. Umweltfreundliche Bauweisen und infrastrukturelle Lösungen in Städten werden zur Grundlage einer widerstandsfähigen, grünen Wirtschaft und städtischer Entwicklung.
'@
$ws.Range("C19").Value = $c19text

$b20text = @'
Entwicklung der Menschheit
'@
$ws.Range("B20").Value = $b20text

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 270
$ws.Rows.Item(6).RowHeight = 409.5
$ws.Rows.Item(7).RowHeight = 409.5
$ws.Rows.Item(8).RowHeight = 409.5
$ws.Rows.Item(12).RowHeight = 409.5
$ws.Rows.Item(18).RowHeight = 409.5
$ws.Rows.Item(19).RowHeight = 409.5
$ws.Rows.Item(20).RowHeight = 45

# --- Selection ---
$ws.Range("A1:C28").Select()
